$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix H2 (setor publico BR / brasileiro) - was incorrectly duplicating G2's old text
$ws.Range("H2").Value = "Aqui é a resenha sobre setor público brasileiro."

# Fix G2 (setor publico ES / capixaba) - add trailing period
$ws.Range("G2").Value = "Aqui é a resenha sobre setor público capixaba."

# Update the selected cell/view state to G2
$ws.Range("G2").Select()
